# "Aggiunto link su cosa e osgi"
#
# 1. Two existing hyperlinks (the mokabyte OSGi article and the Adobe
#    "convert jar to bundle" article) are unwrapped back to plain styled
#    runs (the <w:hyperlink> wrapper is removed, the run + its
#    CollegamentoInternet character style + text stay).
# 2. The last two (empty) paragraphs of the document: the first stays
#    empty, the second gets a new "COS?e OSGi:" heading line, followed
#    by a brand-new paragraph containing the new article link text.
# 3. A new (empty) character style "ListLabel 4" is registered in the
#    style sheet.

$d = $word.ActiveDocument

# --- 1. Unwrap the two hyperlinks, keeping their styled run + text ---
# Locate them by their display text (mokabyte OSGi article / Adobe
# "convert jar to bundle" article) instead of a hard-coded index, then
# remove the <w:hyperlink> wrapper while leaving the styled run + text.
for ($i = $d.Hyperlinks.Count; $i -ge 1; $i--) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.TextToDisplay -eq "http://www.mokabyte.it/2010/07/osgi_fw-3/" -or
        $h.TextToDisplay -eq "https://helpx.adobe.com/it/experience-manager/kb/ConvertAJarIntoOsgiBundle.html") {
        $h.Delete()
    }
}

# --- 2. Turn the trailing empty paragraph into the new heading line and
#        append a fresh paragraph with the new link text ---
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "COS?é OSGi:"
$p.Range.InsertParagraphAfter()
$newP = $d.Paragraphs.Item($d.Paragraphs.Count)
$newP.Range.Text = "http://www.mokabyte.it/2010/02/osgi_fw-1/"

# --- 3. Register the new (empty) character style ---
$s = $d.Styles.Add("ListLabel 4", 2)
$s.QuickStyle = $true
